# Slide 3 - "Qu'est-ce que le HTML ?" intro slide
# The commit slightly modifies this intro slide:
#  1. Title: "Qu'est-ce que le css ?" -> "Qu'est-ce que le CSS ?"
#  2. Acronym line: "Hyper Text Markup Language" -> "Cascaded  Style Sheet"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# --- 1) Title shape ("Rectangle 1") --------------------------------------
$titleShape = $s.Shapes.Item(1)
$titlePara  = $titleShape.TextFrame.TextRange.Paragraphs(1, 1)

# Runs (1-based character offsets), before edit:
#   "Qu'est-ce"  -> 1-9     (untouched)
#   " que le "   -> 10-17
#   "css "       -> 18-21
#   "?"          -> 22
# Edit from the end backwards so earlier offsets stay valid.
$titlePara.Characters(22, 1).Text = ""
$titlePara.Characters(18, 4).Text = ""
$titlePara.Characters(10, 8).Text = " que le CSS ?"

# --- 2) Acronym line in "Text Box 2" --------------------------------------
$acronymShape = $s.Shapes.Item(2)
$acronymPara  = $acronymShape.TextFrame.TextRange.Paragraphs(1, 1)

# Runs (1-based character offsets), before edit:
#   "H"       -> 1       (bold)
#   "yper "   -> 2-6
#   "T"       -> 7       (bold)
#   "ext "    -> 8-11
#   "M"       -> 12      (bold)
#   "arkup "  -> 13-18
#   "L"       -> 19      (bold)
#   "anguage" -> 20-26
# Edit from the end backwards so earlier offsets stay valid; only the
# run text is changed, so existing run formatting (bold / color / font)
# is preserved automatically.
$acronymPara.Characters(20, 7).Text = ""
$acronymPara.Characters(19, 1).Text = ""
$acronymPara.Characters(13, 6).Text = "heet"
$acronymPara.Characters(12, 1).Text = "S"
$acronymPara.Characters(8, 4).Text = "tyle "
$acronymPara.Characters(7, 1).Text = "S"
$acronymPara.Characters(2, 5).Text = "ascaded  "
$acronymPara.Characters(1, 1).Text = "C"
